# Refresh the cryptocurrency symbol list: update the "Price" (column D)
# and "Volume(1h)" (column E) figures for this run's snapshot.
#
# These columns are stored as plain text (not numbers), so the new values
# are assigned with a leading apostrophe to force Excel to keep them as
# text (preserving formatting such as trailing zeros and literal "%" signs)
# instead of auto-converting them into Number/Percentage cells. The style
# is then reset to "Normal" so no extra formatting/number-format is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "306.88"
Set-TextValue "E2" "1.27%"
Set-TextValue "D3" "35.94"
Set-TextValue "E3" "-0.06%"
Set-TextValue "D4" "5.004"
Set-TextValue "E4" "-1.33%"
Set-TextValue "D5" "0.08094"
Set-TextValue "E5" "0.30%"
Set-TextValue "D6" "1.944"
Set-TextValue "E6" "-0.95%"
Set-TextValue "D7" "4.148"
Set-TextValue "E7" "2.41%"
Set-TextValue "D8" "7.881"
Set-TextValue "E8" "1.04%"
Set-TextValue "D9" "0.9325"
Set-TextValue "E9" "0.44%"
Set-TextValue "D10" "0.1252"
Set-TextValue "E10" "-16.95%"
Set-TextValue "D11" "0.1904"
Set-TextValue "E11" "0.42%"
Set-TextValue "D12" "0.09218"
Set-TextValue "E12" "2.43%"
Set-TextValue "D13" "0.03517"
Set-TextValue "E13" "1.80%"
Set-TextValue "D14" "0.09928"
Set-TextValue "E14" "0.82%"
Set-TextValue "D15" "0.001422"
Set-TextValue "E15" "2.18%"
Set-TextValue "D16" "0.006505"
Set-TextValue "E16" "12.55%"
Set-TextValue "E17" "2.12%"
Set-TextValue "D18" "3.108"
Set-TextValue "E18" "4.91%"
Set-TextValue "E19" "-0.12%"
Set-TextValue "E20" "2.39%"
Set-TextValue "D21" "5.160"
Set-TextValue "E21" "2.51%"
Set-TextValue "E22" "5.95%"
Set-TextValue "D23" "0.04407"
Set-TextValue "E23" "-2.04%"
Set-TextValue "D24" "0.001234"
Set-TextValue "E24" "2.38%"
Set-TextValue "D25" "0.004722"
Set-TextValue "E25" "-1.87%"
Set-TextValue "E26" "6.01%"
Set-TextValue "D27" "0.0003132"
Set-TextValue "E27" "3.89%"
Set-TextValue "D39" "0.01962"
Set-TextValue "E39" "4.38%"
Set-TextValue "D40" "0.05201"
Set-TextValue "E40" "8.46%"
Set-TextValue "D41" "0.007569"
Set-TextValue "E41" "3.23%"
Set-TextValue "D42" "0.01018"
Set-TextValue "E42" "-3.86%"
Set-TextValue "D43" "0.1376"
Set-TextValue "E43" "2.42%"
Set-TextValue "E44" "-0.18%"
Set-TextValue "D45" "0.01073"
Set-TextValue "E45" "10.37%"
Set-TextValue "D46" "0.00006372"
Set-TextValue "E46" "2.63%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.45%"
Set-TextValue "D48" "64.96"
Set-TextValue "E48" "0.45%"
Set-TextValue "D49" "0.001661"
Set-TextValue "E49" "0.08%"
Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "0.45%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "0.45%"
